# Add the new "Debug_Messages" setting row to Sheet1 and update the
# active selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New settings row: Debug_Messages flag, default disabled (0)
$ws.Range("A3").Value = "Debug_Messages"
$ws.Range("B3").Value = 0

# Move/restore the visible selection to D4, as recorded in the saved view
[void]$ws.Range("D4").Select()
